$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 10946.1
$ws.Range("J9").Value = 1925
$ws.Range("L9").Value = 1925
$ws.Range("N9").Value = -2263
$ws.Range("H92").Value = 2221.818
$ws.Range("I92").Value = 2249.353
$ws.Range("K92").Value = 2249.353
$ws.Range("M92").Value = -1001.353
$ws.Range("H107").Value = 2037.7142
$ws.Range("I107").Value = 1641
$ws.Range("K107").Value = 1641
$ws.Range("M107").Value = 279
$ws.Range("H118").Value = 2882.1667
$ws.Range("I118").Value = 2723.75
$ws.Range("J118").Value = 3199
$ws.Range("K118").Value = 8171.25
$ws.Range("L118").Value = 9597
$ws.Range("M118").Value = -6514.25
$ws.Range("N118").Value = -12911
$ws.Range("H125").Value = 3738.6
$ws.Range("I125").Value = 2049.25
$ws.Range("J125").Value = 4864.8335
$ws.Range("K125").Value = 18443.25
$ws.Range("L125").Value = 43783.5015
$ws.Range("M125").Value = -15983.25
$ws.Range("N125").Value = -48703.5015
$ws.Range("H138").Value = 5656398
$ws.Range("I138").Value = 3029.2942
$ws.Range("J138").Value = 7944666
$ws.Range("K138").Value = 9087.882599999999
$ws.Range("L138").Value = 23833998
$ws.Range("M138").Value = -3947.882599999999
$ws.Range("N138").Value = -23844278

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 14293367
$ws.Range("I32").Value = 15878218
$ws.Range("J32").Value = 29711.715
$ws.Range("K32").Value = 15878218
$ws.Range("L32").Value = 29711.715
$ws.Range("M32").Value = -15877931
$ws.Range("N32").Value = -30285.715
$ws.Range("H74").Value = 250282530
$ws.Range("I74").Value = 250282530
$ws.Range("K74").Value = 250282530
$ws.Range("M74").Value = -250281656
$ws.Range("H77").Value = 250282530
$ws.Range("I77").Value = 250282530
$ws.Range("K77").Value = 1251412650
$ws.Range("M77").Value = -1251408282
$ws.Range("H88").Value = 2287.4443
$ws.Range("I88").Value = 2430.6667
$ws.Range("K88").Value = 2430.6667
$ws.Range("M88").Value = -2024.6667
$ws.Range("H91").Value = 2287.4443
$ws.Range("I91").Value = 2430.6667
$ws.Range("K91").Value = 2430.6667
$ws.Range("M91").Value = -1026.6667
$ws.Range("H132").Value = 38464748
$ws.Range("I132").Value = 3396.1
$ws.Range("J132").Value = 166669260
$ws.Range("K132").Value = 10188.3
$ws.Range("L132").Value = 500007780
$ws.Range("M132").Value = -7658.299999999999
$ws.Range("N132").Value = -500012840
$ws.Range("H135").Value = 23728.75
$ws.Range("J135").Value = 23728.75
$ws.Range("L135").Value = 23728.75
$ws.Range("N135").Value = -33868.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H81").Value = 56978.332
$ws.Range("J81").Value = 56978.332
$ws.Range("L81").Value = 56978.332
$ws.Range("N81").Value = -59100.332
$ws.Range("H84").Value = 56978.332
$ws.Range("J84").Value = 56978.332
$ws.Range("L84").Value = 170934.996
$ws.Range("N84").Value = -181542.996

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2710.9
$ws.Range("I16").Value = 3099.3333
$ws.Range("K16").Value = 3099.3333
$ws.Range("M16").Value = -2812.3333
$ws.Range("H23").Value = 1000000
$ws.Range("I23").Value = 1000000
$ws.Range("K23").Value = 1000000
$ws.Range("M23").Value = -999760
$ws.Range("H26").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("L26").Value = 0
$ws.Range("N26").ClearContents()
$ws.Range("H27").Value = 1000000
$ws.Range("I27").Value = 1000000
$ws.Range("K27").Value = 1000000
$ws.Range("M27").Value = -999808
$ws.Range("H39").Value = 11666.667
$ws.Range("I39").Value = 11666.667
$ws.Range("K39").Value = 11666.667
$ws.Range("M39").Value = -11275.667
$ws.Range("H49").Value = 11666.667
$ws.Range("I49").Value = 11666.667
$ws.Range("K49").Value = 11666.667
$ws.Range("M49").Value = -11484.667
$ws.Range("H107").Value = 1821.6957
$ws.Range("I107").Value = 1043.6666
$ws.Range("K107").Value = 1043.6666
$ws.Range("M107").Value = 876.3334
$ws.Range("H113").Value = 2710.9
$ws.Range("I113").Value = 3099.3333
$ws.Range("K113").Value = 3099.3333
$ws.Range("M113").Value = -929.3332999999998
$ws.Range("H134").Value = 1089.9744
$ws.Range("I134").Value = 1002.7143
$ws.Range("J134").Value = 1853.5
$ws.Range("K134").Value = 3008.1429
$ws.Range("L134").Value = 5560.5
$ws.Range("M134").Value = -473.1428999999998
$ws.Range("N134").Value = -10630.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 3441.7
$ws.Range("J75").Value = 3839
$ws.Range("L75").Value = 11517
$ws.Range("N75").Value = -13513
$ws.Range("H78").Value = 3441.7
$ws.Range("J78").Value = 3839
$ws.Range("L78").Value = 34551
$ws.Range("N78").Value = -44535
$ws.Range("H134").Value = 4037.4666
$ws.Range("I134").Value = 1658.6923
$ws.Range("K134").Value = 4976.0769
$ws.Range("M134").Value = 93.92309999999998

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1811.1666
$ws.Range("I102").Value = 1693.4286
$ws.Range("K102").Value = 1693.4286
$ws.Range("M102").Value = -71.42859999999996
$ws.Range("H113").Value = 9177.615
$ws.Range("J113").Value = 9210
$ws.Range("L113").Value = 9210
$ws.Range("N113").Value = -13550
$ws.Range("H138").Value = 48000
$ws.Range("I138").Value = 48000
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 48000
$ws.Range("L138").Value = 0
$ws.Range("M138").Value = -42860
$ws.Range("N138").ClearContents()
$ws.Range("H140").Value = 200780
$ws.Range("J140").Value = 200780
$ws.Range("L140").Value = 200780
$ws.Range("N140").Value = -211140

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("M10").ClearContents()
$ws.Range("N10").ClearContents()
$ws.Range("H38").Value = 12400
$ws.Range("J38").Value = 12000
$ws.Range("L38").Value = 12000
$ws.Range("N38").Value = -12820
$ws.Range("H68").Value = 6016.5
$ws.Range("H71").Value = 6016.5
$ws.Range("H82").Value = 3386.7058
$ws.Range("I82").Value = 1175.6666
$ws.Range("K82").Value = 1175.6666
$ws.Range("M82").Value = -814.6666
$ws.Range("H85").Value = 3386.7058
$ws.Range("I85").Value = 1175.6666
$ws.Range("K85").Value = 1175.6666
$ws.Range("M85").Value = 72.33339999999998
$ws.Range("H100").Value = 3382.5151
$ws.Range("I100").Value = 3139.25
$ws.Range("J100").Value = 3611.4707
$ws.Range("K100").Value = 3139.25
$ws.Range("L100").Value = 3611.4707
$ws.Range("M100").Value = -2598.25
$ws.Range("N100").Value = -4693.4707
$ws.Range("H109").Value = 50000
$ws.Range("I109").Value = 50000
$ws.Range("K109").Value = 50000
$ws.Range("M109").Value = -48613
$ws.Range("H124").Value = 46770.6
$ws.Range("J124").Value = 46770.6
$ws.Range("L124").Value = 46770.6
$ws.Range("N124").Value = -56590.6
$ws.Range("H132").Value = 44450384
$ws.Range("I132").Value = 6104.793
$ws.Range("K132").Value = 18314.379
$ws.Range("M132").Value = -15784.379
$ws.Range("H133").Value = 63473.668
$ws.Range("J133").Value = 63473.668
$ws.Range("L133").Value = 63473.668
$ws.Range("N133").Value = -68533.66800000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H10").Value = 20000006
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 20000006
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 20000006
$ws.Range("M10").ClearContents()
$ws.Range("N10").Value = -20000344
$ws.Range("H13").Value = 3100
$ws.Range("J13").Value = 3300
$ws.Range("L13").Value = 3300
$ws.Range("N13").Value = -3580
$ws.Range("H81").Value = 977.2727
$ws.Range("I81").Value = 977.2727
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 1954.5454
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = -893.5454
$ws.Range("N81").ClearContents()
$ws.Range("H84").Value = 977.2727
$ws.Range("I84").Value = 977.2727
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 9772.726999999999
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = -4468.726999999999
$ws.Range("N84").ClearContents()
$ws.Range("H132").Value = 5294.8
$ws.Range("I132").Value = 5315.3613
$ws.Range("J132").Value = 5212.5557
$ws.Range("K132").Value = 15946.0839
$ws.Range("L132").Value = 15637.6671
$ws.Range("M132").Value = -13416.0839
$ws.Range("N132").Value = -20697.6671
$ws.Range("H136").Value = 1712.75
$ws.Range("I136").Value = 1700.8
$ws.Range("J136").Value = 1732.6666
$ws.Range("K136").Value = 5102.4
$ws.Range("L136").Value = 5197.9998
$ws.Range("M136").Value = -2552.4
$ws.Range("N136").Value = -10297.9998
